# Append the latest Pick 3 draw result as a new row at the bottom of the
# "Results" sheet (mirrors the nightly auto-update job that appends one row
# per day: Date, Game, Phase, Result, InsertedAt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row right after the current data (row 65 -> 66).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$target = $ws.Range("A" + $newRow + ":E" + $newRow)

# Force plain text entry so date-looking / digit-only values ("2025-11-21",
# "251121") are stored literally instead of being auto-coerced into a date
# serial / number, matching the rest of the column's text values.
$target.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-11-21"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "251121"
$ws.Cells.Item($newRow, 4).Value = "5-6-6"
$ws.Cells.Item($newRow, 5).Value = "2025-11-21T21:38:49.872+04:00"

# Drop the temporary "@" formatting so the new cells keep the same
# (default/general) style as every other row in the sheet.
$target.ClearFormats()
